$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 576.3387
$ws.Range("J17").Value = 591.7458
$ws.Range("L17").Value = 1775.2374
$ws.Range("N17").Value = -2111.2374
$ws.Range("H98").Value = 83376020
$ws.Range("I98").Value = 83376020
$ws.Range("K98").Value = 83376020
$ws.Range("M98").Value = -83374522
$ws.Range("H99").Value = 491.08334
$ws.Range("I99").Value = 453.9091
$ws.Range("K99").Value = 1361.7273
$ws.Range("M99").Value = 136.2727
$ws.Range("H112").Value = 1332.6666
$ws.Range("H116").Value = 6352.143
$ws.Range("I116").Value = 6100
$ws.Range("J116").Value = 6604.2856
$ws.Range("K116").Value = 6100
$ws.Range("L116").Value = 6604.2856
$ws.Range("M116").Value = -2658
$ws.Range("N116").Value = -13488.2856
$ws.Range("H122").Value = 83376020
$ws.Range("I122").Value = 83376020
$ws.Range("K122").Value = 250128060
$ws.Range("M122").Value = -250125610
$ws.Range("H133").Value = 87999
$ws.Range("J133").Value = 87999
$ws.Range("L133").Value = 87999
$ws.Range("N133").Value = -98119
$ws.Range("H135").Value = 18816.5
$ws.Range("J135").Value = 51000
$ws.Range("L135").Value = 459000
$ws.Range("N135").Value = -464070
$ws.Range("H137").Value = 3578.389
$ws.Range("J137").Value = 7274.25
$ws.Range("L137").Value = 21822.75
$ws.Range("N137").Value = -26922.75
$ws.Range("H138").Value = 1999.375
$ws.Range("I138").Value = 864.05
$ws.Range("J138").Value = 2630.111
$ws.Range("K138").Value = 2592.15
$ws.Range("L138").Value = 7890.333
$ws.Range("M138").Value = 2547.85
$ws.Range("N138").Value = -18170.333
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 2114.1177
$ws.Range("I141").Value = 2183.75
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 6551.25
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = -1371.25
$ws.Range("N141").Value = -13360

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 325
$ws.Range("J5").Value = 466.2857
$ws.Range("L5").Value = 466.2857
$ws.Range("N5").Value = -690.2857
$ws.Range("H32").Value = 17861958
$ws.Range("I32").Value = 19234360
$ws.Range("K32").Value = 19234360
$ws.Range("M32").Value = -19234073
$ws.Range("H45").Value = 3500
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 3500
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -4254

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 325
$ws.Range("J4").Value = 466.2857
$ws.Range("L4").Value = 466.2857
$ws.Range("N4").Value = -696.2857

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1118.8462
$ws.Range("I58").Value = 947.6667
$ws.Range("J58").Value = 1504
$ws.Range("K58").Value = 947.6667
$ws.Range("L58").Value = 1504
$ws.Range("M58").Value = -744.6667
$ws.Range("N58").Value = -1910
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H99").Value = 2953
$ws.Range("J99").Value = 2733.3333
$ws.Range("L99").Value = 2733.3333
$ws.Range("N99").Value = -5729.3333
$ws.Range("H107").Value = 410.53845
$ws.Range("I107").Value = 283.22223
$ws.Range("K107").Value = 283.22223
$ws.Range("M107").Value = 1636.77777
$ws.Range("H122").Value = 1573.5
$ws.Range("I122").Value = 1564.6666
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 4693.9998
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -2243.9998
$ws.Range("N122").Value = -9700
$ws.Range("H126").Value = 2953
$ws.Range("J126").Value = 2733.3333
$ws.Range("L126").Value = 8199.999899999999
$ws.Range("N126").Value = -13139.9999
$ws.Range("H136").Value = 1118.8462
$ws.Range("I136").Value = 947.6667
$ws.Range("J136").Value = 1504
$ws.Range("K136").Value = 2843.0001
$ws.Range("L136").Value = 4512
$ws.Range("M136").Value = -293.0001000000002
$ws.Range("N136").Value = -9612

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9531.875
$ws.Range("I3").Value = 4247.125
$ws.Range("K3").Value = 12741.375
$ws.Range("M3").Value = -12629.375
$ws.Range("H11").Value = 2244.2917
$ws.Range("I11").Value = 2281
$ws.Range("K11").Value = 6843
$ws.Range("M11").Value = -6703
$ws.Range("H117").Value = 974.125
$ws.Range("I117").Value = 1166.5
$ws.Range("J117").Value = 397
$ws.Range("K117").Value = 3499.5
$ws.Range("L117").Value = 1191
$ws.Range("M117").Value = -57.5
$ws.Range("N117").Value = -8075

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1898.8
$ws.Range("I122").Value = 1873.75
$ws.Range("K122").Value = 5621.25
$ws.Range("M122").Value = -3171.25
$ws.Range("H132").Value = 111138220
$ws.Range("I132").Value = 142863420
$ws.Range("K132").Value = 428590260
$ws.Range("M132").Value = -428587730

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4345
$ws.Range("I16").Value = 4517.6665
$ws.Range("K16").Value = 4517.6665
$ws.Range("M16").Value = -4347.6665
$ws.Range("H26").Value = 10009.444
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9705
$ws.Range("H61").Value = 906.9231
$ws.Range("I61").Value = 878.8
$ws.Range("K61").Value = 878.8
$ws.Range("M61").Value = -676.8
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 906.9231
$ws.Range("I113").Value = 878.8
$ws.Range("K113").Value = 878.8
$ws.Range("M113").Value = 1291.2
$ws.Range("H136").Value = 30691.88
$ws.Range("I136").Value = 4164.2666
$ws.Range("J136").Value = 97010.914
$ws.Range("K136").Value = 12492.7998
$ws.Range("L136").Value = 291032.742
$ws.Range("M136").Value = -9942.799800000001
$ws.Range("N136").Value = -296132.742

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 8272.477000000001
$ws.Range("I122").Value = 4223.5
$ws.Range("K122").Value = 12670.5
$ws.Range("M122").Value = -10220.5
$ws.Range("H132").Value = 8616.666999999999
$ws.Range("I132").Value = 1459.2858
$ws.Range("J132").Value = 33667.5
$ws.Range("K132").Value = 4377.857400000001
$ws.Range("L132").Value = 101002.5
$ws.Range("M132").Value = -1847.857400000001
$ws.Range("N132").Value = -106062.5
$ws.Range("H135").Value = 66700
$ws.Range("J135").Value = 67833.336
$ws.Range("L135").Value = 67833.336
$ws.Range("N135").Value = -77973.336
